$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab / workbook.xml <sheet name="..."> from "Evaluation" to "Sheet1"
$ws.Name = "Sheet1"

# Update the Cohere_Translation (column C) text and the recomputed BLEU_Score (column D)
# for the rows whose machine translation text was corrected.

$ws.Range("C17").Value = "सौभाग्यशालीता, ड्रायवर्स आम्हाला सामान्य रूपाने एकच क्षेत्रात ठेवण्याची प्रवृत्ती दाखवतात; त्यामुळे, यातायात धारासामान्य रूपाने काही सामान्य रूपाने एकरूपता असते आणि ती गणितीय रूपाने काही कळण्याची कोशिश केली जाऊ शकते."
$ws.Range("D17").Value = 0.043389519964411

$ws.Range("C30").Value = "दिसावी हजरात वर्षे पूर्वीर, ईरानातील जाग्रोस पर्वतात बकर्यांचा प्रथम घराणीकरण झाला होता."
$ws.Range("D30").Value = 0.01696772192530609

$ws.Range("C57").Value = 'त्याने गॉसिपाला "राजनीतिक गपशप आणि बेवकूफी" म्हणून संबोधित केली.'
$ws.Range("D57").Value = 0.02302676613984314

$ws.Range("C73").Value = "प्रारंभिक रिपोर्ट्स म्हणतात की विमानाला उरुम्चीमध्ये आपात ल्याण्डिंग करायची अनुमती मिळाली नाही यामुळे तो पुन्हा अफ़गानिस्तान कारित कार्यान्वित करण्यात आला."
$ws.Range("D73").Value = 0.01550550784373247

$ws.Range("C78").Value = "सहायक कीट-नाश परीक्षणात भाग घेणारे साफकर्तेयांवर रेंजर्स द्वारे नियंत्रणात असण्यात योजना होती, कारण परीक्षणाची निरीक्षणे केली जात होती आणि त्याची कार्यक्षमता मूल्यांकित केली जात होती."
$ws.Range("D78").Value = 0.01082505792167168

$ws.Range("C80").Value = "त्याच्या प्रवासाच्या दरम्यान, इवासाकी सामुदायिक समस्या सोडविण्यासाठी सामीना आले."
$ws.Range("D80").Value = 0.02573285025273419

$ws.Range("C89").Value = "वॉयेजर्स जो विदेश में लिंग परिवर्तन शल्यचिकित्सा कराना चाहते हैं, उन्हें सुनिश्चित करना चाहिए कि उनके पास वापसी यात्रा के लिए मान्य दस्तावेज हैं।"
$ws.Range("D89").Value = 0.01119961714552871

$ws.Range("C102").Value = "करी हे जड़ीबूटी आणि मसाल्यांवर आधारित डिश आहे ज्यामध्ये मासेच किंवा तरकारीचा समावेश होतो."
$ws.Range("D102").Value = 0.1879831764733509

$ws.Range("C127").Value = "हजारो वर्षे पहिल्या, ग्रीक वैज्ञानिक आरिस्टार्चस म्हणाला की सौर मण्डल सूर्यावर्तुल आहे."
$ws.Range("D127").Value = 0.02004499497906907
